# Recuperação da versão original do Pitch
# Recuperação da versão original do Pitch (1a versão)
#
# The slide had accumulated two extra rectangle shapes ("Retângulo 3",
# id=4 and "Retângulo 4", id=5) that were not part of the original
# (first) version of the Pitch slide. This restores the slide to that
# original state by removing those two shapes, leaving only the title
# ("Título 1", id=2) and subtitle ("Subtítulo 2", id=3) placeholders.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 4 -or $shp.Id -eq 5) {
        $shp.Delete()
    }
}
